$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = 324936160.74
$ws.Range("P2").Value = 236.0770224988
$ws.Range("Q2").Value = 4290675235.13
$ws.Range("R2").Value = 3117.3195119682
$ws.Range("S2").Value = 269794927.13
$ws.Range("T2").Value = 196.0150662736
$ws.Range("U2").Value = -179574977.19
$ws.Range("V2").Value = -130.4672457315
$ws.Range("W2").Value = 100000
$ws.Range("X2").Value = 0.072653355
$ws.Range("Y2").Value = 148009464.97
$ws.Range("Z2").Value = 107.5338420697
$ws.Range("AA2").Value = -7738060.44
$ws.Range("AB2").Value = -5.6219605243
$ws.Range("AC2").Value = 137639892.82
$ws.Range("AD2").Value = 154.5551982345
